# Add new "Greenland" unemployment rate rows (2008-2014) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(2008, 2.0410348667866014),
    @(2009, 3.640364254280859),
    @(2010, 4.5763126850223701),
    @(2011, 5.5478440963271414),
    @(2012, 5.9011802360472094),
    @(2013, 6.072017469584206),
    @(2014, 6.1839003031323676)
)

$row = 32
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = "Greenland"
    $ws.Cells.Item($row, 2).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]
    $ws.Cells.Item($row, 3).NumberFormat = "0.00"
    $row++
}

$ws.Range("E20").Select()
